# Franz_Zeitaufwand.xlsx - "KursleiterKurs und KontaktKurs angefangen"
#
# Adds a new time-tracking entry (row 45: date, activity text, hours) and
# bumps the hours already logged on 2020-01-30 (row 44) from 1 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 44: hours logged that day increase from 1 to 2 -------------------
$ws.Range("C44").Value = 2

# --- Row 45: new entry -----------------------------------------------------
# Copy the date cell's number formatting (style) from the row above so the
# new date cell keeps the same built-in date format instead of creating a
# brand-new style entry.
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A45").Value = 43862   # 2020-02-01
$ws.Range("B45").Value = "Kursbuchung überlegen wie ich das mach, design in visual studio"
$ws.Range("C45").Value = 2

# F3 (=SUM(C4:C100)) and F4 (=180-F3) recalc automatically from the above.

# --- Window / selection state ----------------------------------------------
# Scroll so row 37 is at the top and select B47, matching the author's view
# when they saved.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("B47").Select()
